$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: locate the "1. Team Dashboard" / "a) ...." paragraphs in the
# Appendix section so we can insert the new "Source: ..." reference
# paragraph (plus a following blank paragraph) right before "a) ....".
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $txt = $para.Range.Text
    if ($txt.StartsWith("a) ") -and $txt.Contains("team position")) {
        $target = $para
        break
    }
}

$insertionPoint = $target.Range
$insertionPoint.Collapse(1)
$insertionPoint.InsertParagraphBefore()
$insertionPoint.InsertParagraphBefore()

# Re-resolve the "a) ...." paragraph index after the two new paragraphs
# were spliced in above it.
$aIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $txt = $para.Range.Text
    if ($txt.StartsWith("a) ") -and $txt.Contains("team position")) {
        $aIndex = $i
        break
    }
}

$sourcePara = $d.Paragraphs.Item($aIndex - 2)
$sourceRange = $sourcePara.Range
$sourceRange.Collapse(1)
$sourceStart = $sourceRange.Start

$part1 = "Source: GoodDee. (2020). "
$part2 = "FN312-Investment-Challenge-2020"
$part3 = " [Excel File]. https://github.com/GoodDee/FN312-Investment-Challenge-2020"

$sourceRange.InsertAfter($part1 + $part2 + $part3)

$italicStart = $sourceStart + $part1.Length
$italicEnd = $italicStart + $part2.Length
$italicRange = $d.Range($italicStart, $italicEnd)
$italicRange.Font.Italic = $true
$italicRange.Font.ItalicBi = $true

# ---------------------------------------------------------------------------
# Step 2: relocate the "_GoBack" bookmark from the end of the "d) ...." line
# to the middle of the "c) ...." line (right after the word "historical ").
# ---------------------------------------------------------------------------
$cIndex = $null
$dIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $txt = $para.Range.Text
    if ($txt.StartsWith("c) ") -and $txt.Contains("historical data")) {
        $cIndex = $i
    }
    if ($txt.StartsWith("d) ") -and $txt.Contains("4 line charts")) {
        $dIndex = $i
    }
}

$cPara = $d.Paragraphs.Item($cIndex)
$cFull = $cPara.Range.Text
$splitOffset = $cFull.IndexOf("historical ") + "historical ".Length
$bookmarkPos = $cPara.Range.Start + $splitOffset

$bookmarks = $d.Bookmarks
if ($bookmarks.Exists("_GoBack")) {
    $bookmarks.Item("_GoBack").Delete()
}
$bookmarkTarget = $d.Range($bookmarkPos, $bookmarkPos)
$bookmarks.Add("_GoBack", $bookmarkTarget)

# ---------------------------------------------------------------------------
# Step 3: in the "d) ...." paragraph, move the opening parenthesis from the
# third run into the second run:  "..."  "(4 line charts..." becomes
# "...("  "4 line charts...".
# ---------------------------------------------------------------------------
$dPara = $d.Paragraphs.Item($dIndex)
$dFull = $dPara.Range.Text
$parenOffset = $dFull.IndexOf("(4 line charts")
$parenPos = $dPara.Range.Start + $parenOffset

$parenRange = $d.Range($parenPos, $parenPos + 1)
$parenRange.Delete()
$parenInsert = $d.Range($parenPos, $parenPos)
$parenInsert.InsertAfter("(")

Write-Host "Edits applied"
